$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header cell for new "Dimensions" column, styled bold like the other headers in row 2
$ws.Range("F2").Value = "Dimensions"
$ws.Range("F2").Font.Bold = $true

# Fill F3:F50 with "96X96", right-aligned like the other data cells using that style
$ws.Range("F3:F50").Value = "96X96"
$ws.Range("F3:F50").HorizontalAlignment = -4152

# Size the new column to fit its content (matches the other best-fit columns on this sheet)
$ws.Columns("F").ColumnWidth = 10.666666666666666

# Update the view: drop the frozen/scrolled topLeftCell and select G3 instead
$ws.Range("G3").Select()
